$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.879.97"
$ws.Range("E2").Value = "  -0.22%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.209.87"
$ws.Range("E3").Value = "  -1.64%  "

# Row 4
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.40"
$ws.Range("E5").Value = "  +3.72%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  -0.43%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.88"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.593"
$ws.Range("E9").Value = "  -3.74%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.63"
$ws.Range("E10").Value = "  +1.93%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0901"
$ws.Range("E11").Value = "  -3.71%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.88"
$ws.Range("E12").Value = "  -1.13%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.102"
$ws.Range("E13").Value = "  +0.96%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.537.37"
$ws.Range("E14").Value = "  -1.89%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.40"
$ws.Range("E15").Value = "  -1.42%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.207.33"
$ws.Range("E16").Value = "  -1.74%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.779"
$ws.Range("E17").Value = "  -3.34%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.762.30"
$ws.Range("E18").Value = "  -0.35%  "

# Row 19
$ws.Range("E19").Value = "  -2.71%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.24"
$ws.Range("E20").Value = "  +0.11%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.92"
$ws.Range("E21").Value = "  -1.16%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.61"
$ws.Range("E22").Value = "  -0.61%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.19"
$ws.Range("E23").Value = "  +0.62%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.17"
$ws.Range("E24").Value = "  -8.40%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.63"
$ws.Range("E26").Value = "  -2.09%  "

# Row 27
$ws.Range("B27").Value = "InjectiveProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "39.44"
$ws.Range("E27").Value = "  +1.86%  "

# Row 28
$ws.Range("B28").Value = "WEMIXToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.35"
$ws.Range("E28").Value = "  -2.25%  "

# Row 29
$ws.Range("E29").Value = "  +0.12%  "

# Row 30
$ws.Range("E30").Value = "  -3.39%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.49"
$ws.Range("E31").Value = "  -0.06%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.20"
$ws.Range("E32").Value = "  -0.70%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0844"
$ws.Range("E33").Value = "  +6.04%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.19"
$ws.Range("E34").Value = "  -2.70%  "

# Row 35
$ws.Range("E35").Value = "  -1.29%  "

# Row 36
$ws.Range("E36").Value = "  -3.18%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0351"
$ws.Range("E37").Value = "  +4.79%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.28"
$ws.Range("E38").Value = "  -1.61%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.45"
$ws.Range("E39").Value = "  -3.85%  "

# Row 40
$ws.Range("E40").Value = "  -2.04%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.77"
$ws.Range("E41").Value = "  +15.55%  "

# Row 42
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.198"
$ws.Range("E42").Value = "  -2.95%  "

# Row 43
$ws.Range("B43").Value = "THORChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.26"
$ws.Range("E43").Value = "  -5.31%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "59.54"
$ws.Range("E44").Value = "  -0.99%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.70"
$ws.Range("E45").Value = "  -3.72%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0976"
$ws.Range("E46").Value = "  -1.78%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.26"
$ws.Range("E47").Value = "  -5.21%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.451"
$ws.Range("E48").Value = "  -1.13%  "

# Row 49
$ws.Range("E49").Value = "  -0.39%  "

# Row 50
$ws.Range("E50").Value = "  -1.69%  "

# Row 51
$ws.Range("E51").Value = "  -0.79%  "
